# Update mods data [2025-11-12 07:50:19]
# Append a new row to the ModCounts sheet with the latest mod-count snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (A:C -> 50 chars) -----------------------------------
# ColumnWidth is stored internally with Excel's default 0.8333 char
# padding baked in by this engine, so ask for 50 - 0.8333... to land on
# an exported width of exactly 50.
$ws.Range("A1:C1").ColumnWidth = 49.166666666666664

# --- New row (row 3) ----------------------------------------------------
# Keep the date as literal text (matches the "Date" column already being
# text in row 2) instead of letting Excel coerce it into a date serial.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2025/11/12"
$ws.Range("B3").Value = "逃离鸭科夫"
$ws.Range("C3").Value = 1057

# Center the new row's cells horizontally and vertically.
$newRow = $ws.Range("A3:C3")
$newRow.HorizontalAlignment = -4108   # xlCenter
$newRow.VerticalAlignment = -4108     # xlVAlignCenter
